# Time_recording.xlsx - "added abstract, minutes and requirements"
#
# Fills in the Type-of-Work / Duration / Comment entries for several
# previously-empty activity rows on the "Std-C" sheet. The weekly and
# overview totals are driven entirely by formulas, so they recalculate
# automatically once the raw minutes are entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Std-C")

# --- Week 1 (rows 6-11) ---------------------------------------------------
$ws.Range("B7").Value = "Individual Work"
$ws.Range("C7").Value = 180
$ws.Range("E7").Value = "Brainstorming, looking for ideas"

# --- Week 2 (rows 17-22) ---------------------------------------------------
$ws.Range("E17").Value = "Brainstorming, searching web for ideas"

$ws.Range("B18").Value = "Individual Work"
$ws.Range("C18").Value = 180
$ws.Range("E18").Value = "Implementation"

# --- Week 3 (rows 28-33) ---------------------------------------------------
$ws.Range("B28").Value = "Individual Work"
$ws.Range("C28").Value = 120
$ws.Range("E28").Value = "Brainstorming, forming the actual idea"

$ws.Range("B29").Value = "Individual Work"
$ws.Range("C29").Value = 180
$ws.Range("E29").Value = "Implementation, Documentation"

# --- Week 4 (rows 39-44) ---------------------------------------------------
$ws.Range("B39").Value = "Team-Meeting"
$ws.Range("C39").Value = 150
$ws.Range("E39").Value = "Brainstorming,Implementation, Documentation"

$ws.Range("B40").Value = "Team-Meeting"
$ws.Range("C40").Value = 360
$ws.Range("E40").Value = "Implementation, Documentation"

# --- Week 5 (rows 50-55) ---------------------------------------------------
$ws.Range("B50").Value = "Team-Meeting"
$ws.Range("C50").Value = 240
$ws.Range("E50").Value = "Implementation, Documentation"

$ws.Range("B51").Value = "Individual Work"
$ws.Range("C51").Value = 180
$ws.Range("E51").Value = "Implementation"

# --- Week 6 (rows 61-66) ---------------------------------------------------
$ws.Range("B61").Value = "Team-Meeting"
$ws.Range("C61").Value = 240
$ws.Range("E61").Value = "Implementation, Documentation, Design"

$ws.Range("B62").Value = "Individual Work"
$ws.Range("C62").Value = 480
$ws.Range("E62").Value = "Implementation, Documentation"

$ws.Range("B63").Value = "Individual Work"
$ws.Range("C63").Value = 120
$ws.Range("E63").Value = "Design"

# --- Week 7 (rows 72-77) ---------------------------------------------------
$ws.Range("B72").Value = "Team-Meeting"
$ws.Range("C72").Value = 420
$ws.Range("E72").Value = "Implementation, Documentation"

$ws.Range("B73").Value = "Individual Work"
$ws.Range("C73").Value = 120
$ws.Range("E73").Value = "Implementation"

$ws.Range("B74").Value = "Team-Meeting"
$ws.Range("C74").Value = 360
$ws.Range("E74").Value = "Documentation"

# --- Week 8 (rows 83-88) ---------------------------------------------------
$ws.Range("B83").Value = "Team-Meeting"
$ws.Range("C83").Value = 480
$ws.Range("E83").Value = "Implementation"

$ws.Range("B84").Value = "Individual Work"
$ws.Range("C84").Value = 240
$ws.Range("E84").Value = "Presentation "

$ws.Range("B85").Value = "Individual Work"
$ws.Range("C85").Value = 160
$ws.Range("E85").Value = "Documentation"

# --- Window / selection state ----------------------------------------------
# Std-C becomes the active sheet/tab, with its selection left on E62.
$ws.Activate()
$ws.Range("E62").Select()
